# Adds new "Outbound *" and "Child 1 *" columns to the booking export sheet
# and appends four new booking rows (15-18) pulled from the latest export,
# matching commit "arregado lo del pais".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

function Set-Num($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

# Some numeric-looking values (ages) are stored as text in this export, just
# like the pre-existing rows (e.g. W14). Force a text number format so the
# value round-trips as text instead of being coerced to a number.
function Set-NumAsText($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
}

# ----------------------------------------------------------------------
# New header columns (row 1): AA..AK
# ----------------------------------------------------------------------
Set-Text $ws "AA1" "Outbound Flight Number"
Set-Text $ws "AB1" "Outbound Departure Date"
Set-Text $ws "AC1" "Outbound Departure Time"
Set-Text $ws "AD1" "Outbound Arrival Date"
Set-Text $ws "AE1" "Outbound Arrival Time"
Set-Text $ws "AF1" "Outbound Duration"
Set-Text $ws "AG1" "Child 1 Name"
Set-Text $ws "AH1" "Child 1 Surname"
Set-Text $ws "AI1" "Child 1 Age"
Set-Text $ws "AJ1" "Child 1 Nationality"
Set-Text $ws "AK1" "Child 1 Assistance"

# ----------------------------------------------------------------------
# Row 15 - MEKEC
# ----------------------------------------------------------------------
Set-Text $ws "A15" "MEKEC"
Set-Text $ws "C15" "5/3/2025, 12:41:51"
Set-Text $ws "D15" "SCL"
Set-Text $ws "E15" "BCN"
Set-Text $ws "F15" "One Way"
Set-Text $ws "G15" "Economy"
Set-Text $ws "H15" "Light"
Set-Text $ws "I15" "Economy"
Set-Text $ws "J15" "Light"
Set-Num  $ws "K15" 1
Set-Num  $ws "L15" 0
Set-Num  $ws "M15" 0
Set-Text $ws "N15" "EN"
Set-Text $ws "O15" "Juan"
Set-Text $ws "P15" "Perez"
Set-Text $ws "Q15" "sofiainkoova@gmail.com"
Set-Text $ws "R15" "+93 791234567"
Set-Text $ws "S15" "2222 4000 7000 0005"
Set-Text $ws "T15" "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
Set-Text $ws "U15" "Juan"
Set-Text $ws "V15" "Perez"
Set-NumAsText $ws "W15" "15"
Set-Text $ws "X15" "Afghanistan"
Set-Text $ws "Y15" "Female"
Set-Text $ws "Z15" ""

# ----------------------------------------------------------------------
# Row 16 - MEKH2
# ----------------------------------------------------------------------
Set-Text $ws "A16" "MEKH2"
Set-Text $ws "C16" "5/3/2025, 12:51:51"
Set-Text $ws "D16" "SCL"
Set-Text $ws "E16" "BCN"
Set-Text $ws "F16" "One Way"
Set-Text $ws "G16" "Economy"
Set-Text $ws "H16" "Light"
Set-Text $ws "I16" "Economy"
Set-Text $ws "J16" "Light"
Set-Num  $ws "K16" 1
Set-Num  $ws "L16" 0
Set-Num  $ws "M16" 0
Set-Text $ws "N16" "EN"
Set-Text $ws "O16" "Juan"
Set-Text $ws "P16" "Perez"
Set-Text $ws "Q16" "sofiainkoova@gmail.com"
Set-Text $ws "R16" "+93 791234567"
Set-Text $ws "S16" "2222 4000 7000 0005"
Set-Text $ws "T16" "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
Set-Text $ws "U16" "Juan"
Set-Text $ws "V16" "Perez"
Set-NumAsText $ws "W16" "15"
Set-Text $ws "X16" "Afghanistan"
Set-Text $ws "Y16" "Female"
Set-Text $ws "Z16" ""

# ----------------------------------------------------------------------
# Row 17 - MEKMJ
# ----------------------------------------------------------------------
Set-Text $ws "A17" "MEKMJ"
Set-Text $ws "C17" "5/3/2025, 13:01:48"
Set-Text $ws "D17" "SCL"
Set-Text $ws "E17" "BCN"
Set-Text $ws "F17" "One Way"
Set-Text $ws "G17" "Economy"
Set-Text $ws "H17" "Light"
Set-Text $ws "I17" "Economy"
Set-Text $ws "J17" "Light"
Set-Num  $ws "K17" 1
Set-Num  $ws "L17" 0
Set-Num  $ws "M17" 0
Set-Text $ws "N17" "EN"
Set-Text $ws "O17" "Juan"
Set-Text $ws "P17" "Perez"
Set-Text $ws "Q17" "sofiainkoova@gmail.com"
Set-Text $ws "R17" "+93 791234567"
Set-Text $ws "S17" "2222 4000 7000 0005"
Set-Text $ws "T17" "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
Set-Text $ws "U17" "Juan"
Set-Text $ws "V17" "Perez"
Set-NumAsText $ws "W17" "15"
Set-Text $ws "X17" "Afghanistan"
Set-Text $ws "Y17" "Female"
Set-Text $ws "Z17" ""
Set-Text $ws "AA17" ""
Set-Text $ws "AB17" ""
Set-Text $ws "AC17" ""
Set-Text $ws "AD17" ""
Set-Text $ws "AE17" ""
Set-Text $ws "AF17" ""

# ----------------------------------------------------------------------
# Row 18 - MEKQE
# ----------------------------------------------------------------------
Set-Text $ws "A18" "MEKQE"
Set-Text $ws "C18" "5/3/2025, 13:13:15"
Set-Text $ws "D18" "SCL"
Set-Text $ws "E18" "BCN"
Set-Text $ws "F18" "One Way"
Set-Text $ws "G18" "Economy"
Set-Text $ws "H18" "Light"
Set-Text $ws "I18" "Economy"
Set-Text $ws "J18" "Light"
Set-Num  $ws "K18" 1
Set-Num  $ws "L18" 1
Set-Num  $ws "M18" 0
Set-Text $ws "N18" "EN"
Set-Text $ws "O18" "Juan"
Set-Text $ws "P18" "Perez"
Set-Text $ws "Q18" "sofiainkoova@gmail.com"
Set-Text $ws "R18" "+93 791234567"
Set-Text $ws "S18" "5555 4444 3333 1111"
Set-Text $ws "T18" "Test Consumer"
Set-Text $ws "U18" "Juan"
Set-Text $ws "V18" "Perez"
Set-NumAsText $ws "W18" "15"
Set-Text $ws "X18" "Afghanistan"
Set-Text $ws "Y18" "Female"
Set-Text $ws "Z18" ""
Set-Text $ws "AA18" ""
Set-Text $ws "AB18" ""
Set-Text $ws "AC18" ""
Set-Text $ws "AD18" ""
Set-Text $ws "AE18" ""
Set-Text $ws "AF18" ""
Set-Text $ws "AG18" "Mateo"
Set-Text $ws "AH18" "Ramírez"
Set-NumAsText $ws "AI18" "11"
Set-Text $ws "AJ18" "Argentina"
Set-Text $ws "AK18" "Visual difficulty"
